$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Data correction: vendor names were wrong for two ranges of rows ---
# Rows 52-59 (C column) were labeled "Home Selects" -> should be "HomeSelects"
$ws.Range("C52:C59").Value = "HomeSelects"

# Rows 60-74 (C column) were labeled "Zaca" -> should be "HomeSelects Cabinets"
$ws.Range("C60:C74").Value = "HomeSelects Cabinets"

# --- Reset the view: scroll back to the top and select B2 instead of F71 ---
$win = $excel.ActiveWindow
$win.ScrollRow = 1
$win.ScrollColumn = 1
$ws.Range("B2").Select()
